$wb = $excel.ActiveWorkbook

# ALC row 130: Technically Still Magic / Ophiotauroskin Magitek Codex
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H130").Value = 43780
$ws.Range("J130").Value = 43780
$ws.Range("L130").Value = 43780
$ws.Range("N130").Value = -53820

# ALC row 132: Fast-forwarding Flora / Growth Formula Lambda
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 102322.28
$ws.Range("I132").Value = 116146.18
$ws.Range("J132").Value = 8319.799999999999
$ws.Range("K132").Value = 348438.54
$ws.Range("L132").Value = 24959.4
$ws.Range("M132").Value = -345908.54
$ws.Range("N132").Value = -30019.4

# ALC row 140: Tome for Tradition / Book of Ra'Kaznar
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H140").Value = 41059.445
$ws.Range("J140").Value = 41151.176
$ws.Range("L140").Value = 41151.176
$ws.Range("N140").Value = -51511.176

# ARM row 32: Ingot We Trust / Steel Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5717.449
$ws.Range("I32").Value = 3813.3044
$ws.Range("J32").Value = 9525.739
$ws.Range("K32").Value = 3813.3044
$ws.Range("L32").Value = 9525.739
$ws.Range("M32").Value = -3526.3044
$ws.Range("N32").Value = -10099.739

# ARM row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1966.52
$ws.Range("I61").Value = 1675.591
$ws.Range("K61").Value = 1675.591
$ws.Range("M61").Value = -1463.591

# ARM row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2502.5193
$ws.Range("I132").Value = 1875.2195
$ws.Range("J132").Value = 4840.636
$ws.Range("K132").Value = 5625.6585
$ws.Range("L132").Value = 14521.908
$ws.Range("M132").Value = -3095.6585
$ws.Range("N132").Value = -19581.908

# ARM row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1966.52
$ws.Range("I136").Value = 1675.591
$ws.Range("K136").Value = 5026.772999999999
$ws.Range("M136").Value = -2476.772999999999

# ARM row 137: Odd Instruments / Cobalt Tungsten Alembic
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H137").Value = 40484.332
$ws.Range("J137").Value = 40484.332
$ws.Range("L137").Value = 40484.332
$ws.Range("N137").Value = -50684.332

# BSM row 99: Meddle in Metal / Oroshigane Ingot
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2622.9412
$ws.Range("I99").Value = 1864.3334
$ws.Range("K99").Value = 1864.3334
$ws.Range("M99").Value = -366.3334

# BSM row 107: The Gold Experience / Deepgold Nugget
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1750
$ws.Range("I107").Value = 1500
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 1500
$ws.Range("L107").Value = 2000
$ws.Range("M107").Value = 420
$ws.Range("N107").Value = -5840

# BSM row 134: Ruthenium Supremium / Ruthenium Ingot
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2300.3418
$ws.Range("I134").Value = 1401.0541
$ws.Range("J134").Value = 3092.5715
$ws.Range("K134").Value = 4203.1623
$ws.Range("L134").Value = 9277.7145
$ws.Range("M134").Value = -1668.1623
$ws.Range("N134").Value = -14347.7145

# BSM row 137: Dagger Swagger / Cobalt Tungsten Khukuri
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H137").Value = 38970
$ws.Range("I137").Value = 25000
$ws.Range("J137").Value = 40522.223
$ws.Range("K137").Value = 25000
$ws.Range("L137").Value = 40522.223
$ws.Range("M137").Value = -19900
$ws.Range("N137").Value = -50722.223

# CRP row 31: Wall Not Found / Walnut Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9806555
$ws.Range("I31").Value = 1444.1471
$ws.Range("J31").Value = 29416776
$ws.Range("K31").Value = 1444.1471
$ws.Range("L31").Value = 29416776
$ws.Range("M31").Value = -1149.1471
$ws.Range("N31").Value = -29417366

# CRP row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 9806555
$ws.Range("I34").Value = 1444.1471
$ws.Range("J34").Value = 29416776
$ws.Range("K34").Value = 1444.1471
$ws.Range("L34").Value = 29416776
$ws.Range("M34").Value = -1242.1471
$ws.Range("N34").Value = -29417180

# CRP row 58: You Do the Heavy Lifting / Mahogany Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2217.566
$ws.Range("I58").Value = 1881.0217
$ws.Range("J58").Value = 4429.143
$ws.Range("K58").Value = 1881.0217
$ws.Range("L58").Value = 4429.143
$ws.Range("M58").Value = -1678.0217
$ws.Range("N58").Value = -4835.143

# CRP row 99: O Pine / Pine Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 13337475
$ws.Range("I99").Value = 33336156
$ws.Range("J99").Value = 5022.1113
$ws.Range("K99").Value = 33336156
$ws.Range("L99").Value = 5022.1113
$ws.Range("M99").Value = -33334658
$ws.Range("N99").Value = -8018.1113

# CRP row 107: Built to Last / White Oak Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 977.6
$ws.Range("I107").Value = 630.75
$ws.Range("J107").Value = 1374
$ws.Range("K107").Value = 630.75
$ws.Range("L107").Value = 1374
$ws.Range("M107").Value = 1289.25
$ws.Range("N107").Value = -5214

# CRP row 126: A Better Conductor / Red Pine Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 13337475
$ws.Range("I126").Value = 33336156
$ws.Range("J126").Value = 5022.1113
$ws.Range("K126").Value = 100008468
$ws.Range("L126").Value = 15066.3339
$ws.Range("M126").Value = -100005998
$ws.Range("N126").Value = -20006.3339

# CRP row 132: Hull Lotta Damage / Ginseng Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2831.4412
$ws.Range("I132").Value = 1182.8334
$ws.Range("J132").Value = 4686.125
$ws.Range("K132").Value = 3548.5002
$ws.Range("L132").Value = 14058.375
$ws.Range("M132").Value = -1018.5002
$ws.Range("N132").Value = -19118.375

# CRP row 136: Turali Quality / Dark Mahogany Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2217.566
$ws.Range("I136").Value = 1881.0217
$ws.Range("J136").Value = 4429.143
$ws.Range("K136").Value = 5643.0651
$ws.Range("L136").Value = 13287.429
$ws.Range("M136").Value = -3093.0651
$ws.Range("N136").Value = -18387.429

# CUL row 5: What a Sap / Maple Syrup
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1598.125
$ws.Range("I5").Value = 363.18182
$ws.Range("J5").Value = 2643.077
$ws.Range("K5").Value = 1089.54546
$ws.Range("L5").Value = 7929.231000000001
$ws.Range("M5").Value = -977.54546
$ws.Range("N5").Value = -8153.231000000001

# CUL row 131: The Mountain Steeped / Tsai tou Vounou
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 6850215
$ws.Range("J131").Value = 813.5303
$ws.Range("L131").Value = 2440.5909
$ws.Range("N131").Value = -12520.5909

# CUL row 135: Not-so-secret Ingredient / Royal Maple Syrup
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1598.125
$ws.Range("I135").Value = 363.18182
$ws.Range("J135").Value = 2643.077
$ws.Range("K135").Value = 3268.63638
$ws.Range("L135").Value = 23787.693
$ws.Range("M135").Value = -733.6363799999999
$ws.Range("N135").Value = -28857.693

# GSM row 46: Burning the Midnight Oil / Fire Brand
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 32303.334
$ws.Range("J46").Value = 32303.334
$ws.Range("L46").Value = 32303.334
$ws.Range("N46").Value = -32615.334

# GSM row 107: Whetstones for the Workers / Hard Mudstone Whetstone
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 12346514
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 12346514
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 12346514
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -12350354

# GSM row 122: Awarding Academic Excellence / Ametrine
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2992.88
$ws.Range("I122").Value = 1759.7059
$ws.Range("J122").Value = 5613.375
$ws.Range("K122").Value = 5279.1177
$ws.Range("L122").Value = 16840.125
$ws.Range("M122").Value = -2829.1177
$ws.Range("N122").Value = -21740.125

# GSM row 137: Sew Excited / Cobalt Tungsten Needle
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H137").Value = 48952.777
$ws.Range("J137").Value = 54210
$ws.Range("L137").Value = 54210
$ws.Range("N137").Value = -64410

# LTW row 7: Tan Before the Ban / Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4937.6113
$ws.Range("I7").Value = 2697.5715
$ws.Range("J7").Value = 6363.091
$ws.Range("K7").Value = 2697.5715
$ws.Range("L7").Value = 6363.091
$ws.Range("M7").Value = -2585.5715
$ws.Range("N7").Value = -6587.091

# LTW row 40: Best Served Toad / Toad Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6770.1577
$ws.Range("I40").Value = 5446.952
$ws.Range("J40").Value = 8404.706
$ws.Range("K40").Value = 5446.952
$ws.Range("L40").Value = 8404.706
$ws.Range("M40").Value = -5310.952
$ws.Range("N40").Value = -8676.706

# LTW row 126: Battered Books / Saiga Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 4937.6113
$ws.Range("I126").Value = 2697.5715
$ws.Range("J126").Value = 6363.091
$ws.Range("K126").Value = 8092.7145
$ws.Range("L126").Value = 19089.273
$ws.Range("M126").Value = -5622.7145
$ws.Range("N126").Value = -24029.273

# LTW row 132: Tenets of Tanning / Silver Lobo Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2891.3088
$ws.Range("I132").Value = 1625.2927
$ws.Range("J132").Value = 4813.778
$ws.Range("K132").Value = 4875.8781
$ws.Range("L132").Value = 14441.334
$ws.Range("M132").Value = -2345.8781
$ws.Range("N132").Value = -19501.334

# WVR row 107: Flax Wax / Bright Linen Yarn
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 675.7778
$ws.Range("I107").Value = 544.06665
$ws.Range("K107").Value = 1632.19995
$ws.Range("M107").Value = 287.8000500000001

# WVR row 132: Comfy Cabins / Snow Cotton Cloth
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 11113074
$ws.Range("I132").Value = 993.05884
$ws.Range("J132").Value = 25644258
$ws.Range("K132").Value = 2979.17652
$ws.Range("L132").Value = 76932774
$ws.Range("M132").Value = -449.17652
$ws.Range("N132").Value = -76937834

# WVR row 136: Weaving the Envelope / Sarcenet Cloth
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1553
$ws.Range("I136").Value = 660.0714
$ws.Range("J136").Value = 2514.6155
$ws.Range("K136").Value = 1980.2142
$ws.Range("L136").Value = 7543.8465
$ws.Range("M136").Value = 569.7857999999999
$ws.Range("N136").Value = -12643.8465
